# "Elimna EC anteriores y se agregan nuevos, se modifica base de datos"
#
# The worker data table (rows 16-23 of Hoja1) previously listed the four
# workers grouped by "Periodo Mora" (all 1712 rows first, then all 1801
# rows). The database was reorganised so the rows are grouped by worker
# instead, each worker now followed immediately by their two periods
# (1801 then 1712).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Doc = "23071486";   Nombre = "DANILZA MARIA MERIÑO DIAZ";  Periodo = "1801" },
    @{ Doc = "23071486";   Nombre = "DANILZA MARIA MERIÑO DIAZ";  Periodo = "1712" },
    @{ Doc = "32941250";   Nombre = "LUCILA HERRERA GUZMAN";      Periodo = "1801" },
    @{ Doc = "32941250";   Nombre = "LUCILA HERRERA GUZMAN";      Periodo = "1712" },
    @{ Doc = "1049533950"; Nombre = "ROSANA MARIA PADILLA JULIO"; Periodo = "1801" },
    @{ Doc = "1049533950"; Nombre = "ROSANA MARIA PADILLA JULIO"; Periodo = "1712" },
    @{ Doc = "1049532082"; Nombre = "LORENA MARIA BOLAÑOS UTRIA"; Periodo = "1801" },
    @{ Doc = "1049532082"; Nombre = "LORENA MARIA BOLAÑOS UTRIA"; Periodo = "1712" }
)

$startRow = 16
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 3).Value = $rows[$i].Doc
    $ws.Cells.Item($r, 4).Value = $rows[$i].Nombre
    $ws.Cells.Item($r, 5).Value = $rows[$i].Periodo
}
